$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.02185066666666667
$ws.Cells.Item(2, 8).Value = 0.065552
$ws.Cells.Item(2, 9).Value = 0.02597345993572409
$ws.Cells.Item(2, 10).Value = 0.02597345993572409
$ws.Cells.Item(2, 13).Value = 19.741419
$ws.Cells.Item(2, 14).Value = 59.224257
$ws.Cells.Item(2, 15).Value = 0.5456357702458839
$ws.Cells.Item(2, 16).Value = 0.5456357702458838
$ws.Cells.Item(2, 17).Value = 0.431363166096
$ws.Cells.Item(2, 18).Value = 3.882268494864
$ws.Cells.Item(2, 19).Value = 0.01417204881797942
$ws.Cells.Item(2, 20).Value = 0.01417204881797942

$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.02185066666666667
$ws.Cells.Item(3, 8).Value = 0.065552
$ws.Cells.Item(3, 9).Value = 0.02597345993572409
$ws.Cells.Item(3, 10).Value = 0.02597345993572409
$ws.Cells.Item(3, 15).Value = 0.3560475839792585
$ws.Cells.Item(3, 16).Value = 0.3560475839792585
$ws.Cells.Item(3, 17).Value = 0.2814804700888889
$ws.Cells.Item(3, 18).Value = 2.5333242308
$ws.Cells.Item(3, 19).Value = 0.009247787657696627
$ws.Cells.Item(3, 20).Value = 0.009247787657696627

$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.02185066666666667
$ws.Cells.Item(4, 8).Value = 0.065552
$ws.Cells.Item(4, 9).Value = 0.02597345993572409
$ws.Cells.Item(4, 10).Value = 0.02597345993572409
$ws.Cells.Item(4, 15).Value = 0.09831664577485776
$ws.Cells.Item(4, 16).Value = 0.09831664577485774
$ws.Cells.Item(4, 17).Value = 0.07772617176888889
$ws.Cells.Item(4, 18).Value = 0.69953554592
$ws.Cells.Item(4, 19).Value = 0.002553623460048045
$ws.Cells.Item(4, 20).Value = 0.002553623460048044

$ws.Cells.Item(5, 9).Value = 0.6906391812052189
$ws.Cells.Item(5, 10).Value = 0.6906391812052189
$ws.Cells.Item(5, 13).Value = 19.741419
$ws.Cells.Item(5, 14).Value = 59.224257
$ws.Cells.Item(5, 15).Value = 0.5456357702458839
$ws.Cells.Item(5, 16).Value = 0.5456357702458838
$ws.Cells.Item(5, 17).Value = 11.47002765792
$ws.Cells.Item(5, 18).Value = 103.23024892128
$ws.Cells.Item(5, 19).Value = 0.3768374415988963
$ws.Cells.Item(5, 20).Value = 0.3768374415988962

$ws.Cells.Item(6, 9).Value = 0.6906391812052189
$ws.Cells.Item(6, 10).Value = 0.6906391812052189
$ws.Cells.Item(6, 15).Value = 0.3560475839792585
$ws.Cells.Item(6, 16).Value = 0.3560475839792585
$ws.Cells.Item(6, 19).Value = 0.2459004118695315
$ws.Cells.Item(6, 20).Value = 0.2459004118695315

$ws.Cells.Item(7, 9).Value = 0.6906391812052189
$ws.Cells.Item(7, 10).Value = 0.6906391812052189
$ws.Cells.Item(7, 15).Value = 0.09831664577485776
$ws.Cells.Item(7, 16).Value = 0.09831664577485774
$ws.Cells.Item(7, 19).Value = 0.06790132773679131
$ws.Cells.Item(7, 20).Value = 0.06790132773679131

$ws.Cells.Item(8, 9).Value = 0.283387358859057
$ws.Cells.Item(8, 10).Value = 0.283387358859057
$ws.Cells.Item(8, 13).Value = 19.741419
$ws.Cells.Item(8, 14).Value = 59.224257
$ws.Cells.Item(8, 15).Value = 0.5456357702458839
$ws.Cells.Item(8, 16).Value = 0.5456357702458838
$ws.Cells.Item(8, 17).Value = 4.706452996695
$ws.Cells.Item(8, 18).Value = 42.358076970255
$ws.Cells.Item(8, 19).Value = 0.1546262798290083
$ws.Cells.Item(8, 20).Value = 0.1546262798290082

$ws.Cells.Item(9, 9).Value = 0.283387358859057
$ws.Cells.Item(9, 10).Value = 0.283387358859057
$ws.Cells.Item(9, 15).Value = 0.3560475839792585
$ws.Cells.Item(9, 16).Value = 0.3560475839792585
$ws.Cells.Item(9, 19).Value = 0.1008993844520303
$ws.Cells.Item(9, 20).Value = 0.1008993844520303

$ws.Cells.Item(10, 9).Value = 0.283387358859057
$ws.Cells.Item(10, 10).Value = 0.283387358859057
$ws.Cells.Item(10, 15).Value = 0.09831664577485776
$ws.Cells.Item(10, 16).Value = 0.09831664577485774
$ws.Cells.Item(10, 19).Value = 0.0278616945780184
$ws.Cells.Item(10, 20).Value = 0.0278616945780184
